$d = $word.ActiveDocument

# --- 1. Wrap the first "DONE" run (before "The quality, readability...")
#     with a brand-new bookmark __DdeLink__95_1326230006. Adding this
#     bookmark makes Word renumber every bookmark by document order, which
#     is exactly why the two pre-existing __DdeLink__ bookmarks below shift
#     from id=1/2 to id=2/3 -- no manual renumbering needed.
$pQuality = $d.Paragraphs.Item(29)
$qRange = $pQuality.Range
$doneRange = $d.Range($qRange.Start, $qRange.End)
$null = $doneRange.Find.Execute("DONE")
$d.Bookmarks.Add("__DdeLink__95_1326230006", $doneRange)

# --- 2. Prefix the "Live/dead warm classification..." bullet with its own
#     green "DONE" run (matching the styling used by the sibling bullets).
$pLiveDead = $d.Paragraphs.Item(35)
$liveDeadRange = $pLiveDead.Range
$insertPoint = $d.Range($liveDeadRange.Start, $liveDeadRange.Start)
$insertPoint.InsertBefore("DONE")
$doneRun = $d.Range($liveDeadRange.Start, $liveDeadRange.Start + 4)
$doneRun.Font.Name = "Arial"
$doneRun.Font.Color = 52224

# --- 3. Append the 54 new ListLabel185..ListLabel238 character styles right
#     after the existing ListLabel184 style definition.
$xml = $d.WordOpenXML
$marker = '<w:style w:type="character" w:styleId="ListLabel184"><w:name w:val="ListLabel 184"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style>'
$idx = $xml.IndexOf($marker)
if ($idx -lt 0) {
    throw "ListLabel184 style marker not found in WordOpenXML"
}
$insertAt = $idx + $marker.Length
$newStyles = '<w:style w:type="character" w:styleId="ListLabel185"><w:name w:val="ListLabel 185"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/><w:sz w:val="24"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel186"><w:name w:val="ListLabel 186"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel187"><w:name w:val="ListLabel 187"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel188"><w:name w:val="ListLabel 188"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel189"><w:name w:val="ListLabel 189"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel190"><w:name w:val="ListLabel 190"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel191"><w:name w:val="ListLabel 191"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel192"><w:name w:val="ListLabel 192"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel193"><w:name w:val="ListLabel 193"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel194"><w:name w:val="ListLabel 194"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel195"><w:name w:val="ListLabel 195"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel196"><w:name w:val="ListLabel 196"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel197"><w:name w:val="ListLabel 197"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel198"><w:name w:val="ListLabel 198"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel199"><w:name w:val="ListLabel 199"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel200"><w:name w:val="ListLabel 200"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel201"><w:name w:val="ListLabel 201"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel202"><w:name w:val="ListLabel 202"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel203"><w:name w:val="ListLabel 203"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel204"><w:name w:val="ListLabel 204"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel205"><w:name w:val="ListLabel 205"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel206"><w:name w:val="ListLabel 206"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel207"><w:name w:val="ListLabel 207"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel208"><w:name w:val="ListLabel 208"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel209"><w:name w:val="ListLabel 209"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel210"><w:name w:val="ListLabel 210"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel211"><w:name w:val="ListLabel 211"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel212"><w:name w:val="ListLabel 212"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel213"><w:name w:val="ListLabel 213"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel214"><w:name w:val="ListLabel 214"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel215"><w:name w:val="ListLabel 215"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel216"><w:name w:val="ListLabel 216"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel217"><w:name w:val="ListLabel 217"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel218"><w:name w:val="ListLabel 218"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel219"><w:name w:val="ListLabel 219"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel220"><w:name w:val="ListLabel 220"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel221"><w:name w:val="ListLabel 221"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel222"><w:name w:val="ListLabel 222"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel223"><w:name w:val="ListLabel 223"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel224"><w:name w:val="ListLabel 224"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel225"><w:name w:val="ListLabel 225"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel226"><w:name w:val="ListLabel 226"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel227"><w:name w:val="ListLabel 227"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel228"><w:name w:val="ListLabel 228"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel229"><w:name w:val="ListLabel 229"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel230"><w:name w:val="ListLabel 230"/><w:qFormat/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel231"><w:name w:val="ListLabel 231"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel232"><w:name w:val="ListLabel 232"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel233"><w:name w:val="ListLabel 233"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel234"><w:name w:val="ListLabel 234"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel235"><w:name w:val="ListLabel 235"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel236"><w:name w:val="ListLabel 236"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel237"><w:name w:val="ListLabel 237"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style><w:style w:type="character" w:styleId="ListLabel238"><w:name w:val="ListLabel 238"/><w:qFormat/><w:rPr><w:rFonts w:cs="OpenSymbol"/></w:rPr></w:style>'
$xml = $xml.Substring(0, $insertAt) + $newStyles + $xml.Substring($insertAt)
$d.WordOpenXML = $xml

Write-Output "edit complete"
